$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row: "_old" -> "_FV2404", "_new" -> "_FV2410" ---
# Columns A-J (1-10) hold the "_old" headers
$headersOld = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# Column K (11) is "diff" - unchanged

# Columns L-U (12-21) hold the "_new" headers
$headersNew = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headersOld.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersOld[$i]
}

for ($i = 0; $i -lt $headersNew.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersNew[$i]
}

# --- Freeze the header row (pane split at row 1) ---
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Convert the data range into an Excel Table (ListObject) ---
$range = $ws.Range("A1:U74")
$listObject = $ws.ListObjects.Add(1, $range, [System.Type]::Missing, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""
